$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "51.627.62"

Set-TextValue "D3" "3.019.31"
$ws.Range("E3").Value = "  +2.04%  "

Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.08%  "

Set-TextValue "D5" "379.26"
$ws.Range("E5").Value = "  -0.04%  "

Set-TextValue "D6" "102.61"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("E10").Value = "  +0.79%  "

$ws.Range("E11").Value = "  -0.21%  "

Set-TextValue "D12" "0.0862"
$ws.Range("E12").Value = "  +1.29%  "

Set-TextValue "D13" "3.502.34"
$ws.Range("E13").Value = "  +2.07%  "

Set-TextValue "D14" "18.46"
$ws.Range("E14").Value = "  +0.17%  "

Set-TextValue "D16" "3.017.25"
$ws.Range("E16").Value = "  +2.14%  "

$ws.Range("E17").Value = "  -3.70%  "

Set-TextValue "D18" "10.59"
$ws.Range("E18").Value = "  -14.43%  "

Set-TextValue "D19" "51.628.56"
$ws.Range("E19").Value = "  +1.01%  "

$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("E21").Value = "  +0.21%  "

Set-TextValue "D22" "0.0₃0964"
$ws.Range("E22").Value = "  +0.82%  "

$ws.Range("E23").Value = "  +0.50%  "

Set-TextValue "D24" "267.38"
$ws.Range("E24").Value = "  -0.96%  "

$ws.Range("E25").Value = "  -5.95%  "

Set-TextValue "D26" "8.30"
$ws.Range("E26").Value = "  +3.66%  "

$ws.Range("E27").Value = "  +7.48%  "

$ws.Range("E28").Value = "  +4.16%  "

$ws.Range("E29").Value = "  +0.01%  "

Set-TextValue "D30" "26.20"
$ws.Range("E30").Value = "  +1.43%  "

$ws.Range("E31").Value = "  +0.46%  "

Set-TextValue "D32" "10.25"
$ws.Range("E32").Value = "  -2.55%  "

Set-TextValue "D33" "2.11"
$ws.Range("E33").Value = "  +2.35%  "

Set-TextValue "D34" "50.55"
$ws.Range("E34").Value = "  -0.73%  "

Set-TextValue "D35" "33.79"
$ws.Range("E35").Value = "  -1.66%  "

Set-TextValue "D36" "0.0450"
$ws.Range("E36").Value = "  +3.36%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("E38").Value = "  +2.23%  "

$ws.Range("E39").Value = "  +13.40%  "

Set-TextValue "D40" "16.91"
$ws.Range("E40").Value = "  +0.78%  "

$ws.Range("E41").Value = "  +1.31%  "

Set-TextValue "D42" "127.41"
$ws.Range("E42").Value = "  +5.88%  "

$ws.Range("E43").Value = "  -0.73%  "

$ws.Range("E44").Value = "  +1.37%  "

Set-TextValue "D45" "3.78"
$ws.Range("E45").Value = "  +5.21%  "

Set-TextValue "D46" "21.54"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("E47").Value = "  +2.64%  "

Set-TextValue "D48" "2.41"
$ws.Range("E48").Value = "  +2.63%  "

Set-TextValue "D49" "2.025.28"
$ws.Range("E49").Value = "  -1.02%  "

Set-TextValue "D50" "3.318.62"
$ws.Range("E50").Value = "  +1.99%  "

$ws.Range("E51").Value = "  -1.86%  "
